# Adds the listed control function tests to the "Control functions" sheet.
# Several rows that previously had no test results (Tested = "N") now show
# the test was performed (Tested = "Y") along with Pass/Fail, script name
# and remarks; a couple of existing remarks are also refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Control functions")

# --- ControlChooseIndex / ControlChooseString: mark Pass/Fail as Pass ---
$ws.Range("E3").Value = "Pass"
$ws.Range("E4").Value = "Pass"

# --- ControlClick: remarks updated from stale note to WIP ---
$ws.Range("G5").Value = "WIP"

# --- ControlGetFocus: now tested ---
$ws.Range("C13").Value = "Y"
$ws.Range("E13").Value = "Pass"
$ws.Range("F13").Value = "control_get_focus.ahk"

# --- ControlGetIndex: now tested ---
$ws.Range("C15").Value = "Y"
$ws.Range("E15").Value = "Pass"
$ws.Range("F15").Value = "control_get_index.ahk"

# --- ControlGetItems: now tested ---
$ws.Range("C16").Value = "Y"
$ws.Range("E16").Value = "Pass"
$ws.Range("F16").Value = "control_get_items.ahk"
$ws.Range("G16").Value = "Also tested w/combobox"

# --- ControlGetPos: now tested ---
$ws.Range("C17").Value = "Y"
$ws.Range("E17").Value = "Fail"
$ws.Range("G17").Value = "Some functionality is there."

# --- ControlGetExStyle: now tested (also gains the Function description) ---
$ws.Range("B19").Value = "Returns an integer representing the style or extended style of the specified control."
$ws.Range("C19").Value = "Y"
$ws.Range("E19").Value = "Pass"
$ws.Range("F19").Value = "control_exstyle_cb.ahk"
$ws.Range("G19").Value = "Not in guitest.ahk yet"

# --- ControlHideDropDown: now passes ---
$ws.Range("E23").Value = "Pass"

# --- ControlMove: now tested ---
$ws.Range("C24").Value = "Y"
$ws.Range("E24").Value = "Pass"
$ws.Range("F24").Value = "control_move.ahk"

# --- ControlSend: now tested ---
$ws.Range("C25").Value = "Y"
$ws.Range("E25").Value = "Pass"
$ws.Range("F25").Value = "control_send.ahk"

# --- ControlSendText: now tested (also gains the Function description) ---
$ws.Range("B26").Value = "Sends simulated keystrokes or text to a window or control."
$ws.Range("C26").Value = "Y"
$ws.Range("E26").Value = "Pass"
$ws.Range("F26").Value = "control_send.ahk"

# --- ControlSetChecked: now tested ---
$ws.Range("C27").Value = "Y"
$ws.Range("E27").Value = "Fail"
$ws.Range("F27").Value = "control_set_checked.ahk"

# --- ControlSetEnabled: now tested ---
$ws.Range("C28").Value = "Y"
$ws.Range("E28").Value = "Fail"
$ws.Range("F28").Value = "control_set_enabled.ahk"

# --- ControlSetExStyle: now tested (also gains the Function description) ---
$ws.Range("B30").Value = "Changes the style or extended style of the specified control, respectively."
$ws.Range("C30").Value = "Y"
$ws.Range("E30").Value = "Pass"
$ws.Range("F30").Value = "control_exstyle_cb.ahk"

# --- ControlShowDropDown: now passes ---
$ws.Range("E33").Value = "Pass"

# --- ListViewGetContent: now tested under ControlZoo - Group One, passes ---
$ws.Range("D40").Value = "ControlZoo - Group One"
$ws.Range("E40").Value = "Pass"
$ws.Range("G40").ClearContents()

# Restore the active selection to match the authored workbook.
$ws.Range("E33").Select()
